$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting already used in column B (style applied to B48) down
# into the two new rows before writing values, so the new cells reuse the
# existing date style instead of Excel minting a brand-new numFmt.
$ws.Range("B48").Copy() | Out-Null
$ws.Range("B49:B50").PasteSpecial(-4122) | Out-Null

# Row 49: Federico Speroni, 2017-05-24, 2 horas, Sprint 3 - Integración BackEnd y FrontEnd, Investigación sobre API REST
$ws.Cells.Item(49, 1).Value = "Federico Speroni"
$ws.Cells.Item(49, 2).Value = 42879
$ws.Cells.Item(49, 3).Value = 2
$ws.Cells.Item(49, 4).Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Cells.Item(49, 5).Value = "Investigación sobre API REST"

# Row 50: Federico Speroni, 2017-05-25, 2 horas, Sprint 3 - Integración BackEnd y FrontEnd, Investigación sobre API REST
$ws.Cells.Item(50, 1).Value = "Federico Speroni"
$ws.Cells.Item(50, 2).Value = 42880
$ws.Cells.Item(50, 3).Value = 2
$ws.Cells.Item(50, 4).Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Cells.Item(50, 5).Value = "Investigación sobre API REST"

# Update selection to match the diff (B51) and scroll position remains A28
$ws.Range("B51").Select()
